# Generate responses to test questions
# Renames the header row to the lower_snake_case column names used by the
# downstream RAG pipeline and fixes a typo in one of the edited questions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rename (A1:F1) ---
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "question_edited"
$ws.Range("C1").Value = "answer"
$ws.Range("D1").Value = "response_no_context"
$ws.Range("E1").Value = "response_context"
$ws.Range("F1").Value = "source"

# --- Fix typo: "gfts" -> "gifts" in the edited question for row 11 ---
$ws.Range("B11").Value = "Can Government employees accept monetary gifts and entertainment from someone seeking to obtain Government business as long as they are not solicited?"

# --- Update the view: scroll down and move the active selection to B11 ---
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("B11").Select()
